$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -4
$ws.Range("F9").Value = 2
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = -4
$ws.Range("F24").Value = -5
$ws.Range("F26").Value = 5
$ws.Range("F27").Value = -1
$ws.Range("F28").Value = 1
$ws.Range("F30").Value = 6
$ws.Range("F32").Value = 0
$ws.Range("F38").Value = 4
$ws.Range("F42").Value = -3
$ws.Range("F43").Value = -2
$ws.Range("F47").Value = -1
$ws.Range("F48").Value = -5
$ws.Range("F50").Value = 2
$ws.Range("F54").Value = 0
$ws.Range("F55").Value = 5
$ws.Range("F59").Value = 2
$ws.Range("F63").Value = -11
$ws.Range("F64").Value = -1
$ws.Range("F65").Value = 1
$ws.Range("F67").Value = -3
